# Weekly NYPD CompStat data refresh: "New crime data collected"
# Updates the report header (issue number + week-covering dates) and the
# weekly crime-complaint figures (rows 14-31, 33) on the CompStat sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text: volume/number and the week-covering date range.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value  = "Volume 32   Number  8"
$ws.Cells.Item(9, 3).Value  = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# ---------------------------------------------------------------------------
# 2. Weekly crime-complaint figures, columns C..N for rows 14-31 and 33.
#    $null entries mean "leave the existing value/style untouched".
# ---------------------------------------------------------------------------
$rows = @(
    @{Row=14; C=5;    D=6;    E=-16.666666666666; F=21;   G=24;   H=-12.5;            I=43;   J=55;   K=-21.818181818181; L=-25.862068965517; M=-36.764705882352; N=-85.570469798657}
    @{Row=15; C=36;   D=31;   E=16.129032258064;  F=140;  G=114;  H=22.807017543859;  I=267;  J=201;  K=32.835820895522;  L=16.593886462882;  M=64.814814814814;  N=-39.179954441913}
    @{Row=16; C=225;  D=300;  E=-25;              F=942;  G=1212; H=-22.277227722772; I=1834; J=2414; K=-24.026512013256; L=-21.118279569892; M=-31.387953610175; N=-85.496243574535}
    @{Row=17; C=430;  D=453;  E=-5.077262693156;  F=1847; G=1997; H=-7.511266900350;  I=3533; J=3745; K=-5.660881174899;  L=-3.942359978249;  M=63.792304126101;  N=-30.246791707798}
    @{Row=18; C=232;  D=253;  E=-8.300395256917;  F=911;  G=963;  H=-5.399792315680;  I=1808; J=1911; K=-5.389848246991;  L=-19.141323792486; M=-33.869787856620; N=-87.899879534198}
    @{Row=19; C=766;  D=872;  E=-12.155963302752; F=3016; G=3727; H=-19.077005634558; I=5739; J=7167; K=-19.924654667224; L=-18.607289746135; M=13.666072489601;  N=-49.360275302214}
    @{Row=20; C=204;  D=217;  E=-5.990783410138;  F=820;  G=959;  H=-14.494264859228; I=1571; J=1954; K=-19.600818833162; L=-27.201112140871; M=13.923132704858;  N=-90.848721384050}
    @{Row=21; C=1898; D=2132; E=-10.975609756097; F=7697; G=8996; H=-14.439751000444; I=14795;J=17447;K=-15.200320972086; L=-16.577389343106; M=4.028969202643;   N=-76.094297855838}
    @{Row=22; C=33;   D=40;   E=-17.5;            F=141;  G=168;  H=-16.071428571428; I=264;  J=359;  K=-26.462395543175; L=-11.705685618729; M=-15.112540192926; N=$null}
    @{Row=23; C=93;   D=101;  E=-7.920792079207;  F=408;  G=423;  H=-3.546099290780;  I=755;  J=854;  K=-11.592505854800; L=-13.118527042577; M=40.074211502782;  N=$null}
    @{Row=24; C=1997; D=2044; E=-2.299412915851;  F=7966; G=8358; H=-4.690117252931;  I=15038;J=15650;K=-3.910543130990;  L=-4.204357242960;  M=40.647212869435;  N=$null}
    @{Row=25; C=1158; D=1132; E=2.296819787985;   F=4380; G=4661; H=-6.028749195451;  I=8157; J=8727; K=-6.531454107940;  L=-4.830241512075;  M=$null;            N=$null}
    @{Row=26; C=689;  D=776;  E=-11.211340206185; F=2960; G=3318; H=-10.789632308619; I=5807; J=6108; K=-4.927963326784;  L=0.903562119895;   M=-1.742808798646;  N=$null}
    @{Row=27; C=42;   D=44;   E=-4.545454545454;  F=173;  G=174;  H=-0.574712643678;  I=337;  J=337;  K=0;                L=-7.923497267759;  M=$null;            N=$null}
    @{Row=28; C=92;   D=$null;E=13.580246913580;  F=390;  G=356;  H=9.550561797752;   I=713;  J=657;  K=8.523592085235;   L=3.333333333333;   M=$null;            N=$null}
    @{Row=29; C=14;   D=11;   E=27.272727272727;  F=53;   G=65;   H=-18.461538461538; I=105;  J=128;  K=-17.96875;        L=-32.258064516129; M=-42.622950819672; N=-87.241798298906}
    @{Row=30; C=13;   D=8;    E=62.5;             F=51;   G=52;   H=-1.923076923076;  I=86;   J=105;  K=-18.095238095238; L=-32.8125;          M=-44.516129032258; N=-88.713910761154}
    @{Row=31; C=$null;D=5;    E=-100;             F=24;   G=35;   H=-31.428571428571; I=48;   J=66;   K=-27.272727272727; L=-23.809523809523; M=$null;            N=$null}
    @{Row=33; C=$null;D=6;    E=-100;             F=12;   G=16;   H=-25;              I=23;   J=34;   K=-32.352941176470; L=-17.857142857142; M=$null;            N=$null}
)

foreach ($r in $rows) {
    if ($null -ne $r.C) { $ws.Cells.Item($r.Row, 3).Value  = $r.C }
    if ($null -ne $r.D) { $ws.Cells.Item($r.Row, 4).Value  = $r.D }
    if ($null -ne $r.E) { $ws.Cells.Item($r.Row, 5).Value  = $r.E }
    if ($null -ne $r.F) { $ws.Cells.Item($r.Row, 6).Value  = $r.F }
    if ($null -ne $r.G) { $ws.Cells.Item($r.Row, 7).Value  = $r.G }
    if ($null -ne $r.H) { $ws.Cells.Item($r.Row, 8).Value  = $r.H }
    if ($null -ne $r.I) { $ws.Cells.Item($r.Row, 9).Value  = $r.I }
    if ($null -ne $r.J) { $ws.Cells.Item($r.Row, 10).Value = $r.J }
    if ($null -ne $r.K) { $ws.Cells.Item($r.Row, 11).Value = $r.K }
    if ($null -ne $r.L) { $ws.Cells.Item($r.Row, 12).Value = $r.L }
    if ($null -ne $r.M) { $ws.Cells.Item($r.Row, 13).Value = $r.M }
    if ($null -ne $r.N) { $ws.Cells.Item($r.Row, 14).Value = $r.N }
}

# ---------------------------------------------------------------------------
# 3. Rows 31 and 33, column C: crime count dropped to zero -> Excel (and the
#    source feed) renders a week-to-date zero as the literal text "0" rather
#    than a number, matching the formatting already used by the "***.*" /
#    not-applicable cells on these rows. Force text storage, then restore the
#    plain (non-bold, general-format) look of the surrounding text cells.
# ---------------------------------------------------------------------------
$ws.Range("C31").Value = "'0"
$ws.Range("M31").Copy()
$ws.Range("C31").PasteSpecial(-4122)

$ws.Range("C33").Value = "'0"
$ws.Range("M33").Copy()
$ws.Range("C33").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0
